# Insert a new data row at row 218 (shifts existing rows 218-312 down to 219-313)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with its values.
# Columns A,B,C,E,F,G,H,I,J,K keep the same constant values used throughout
# this block of rows (Mercado/Producto metadata); only D,L,M,N,O,P,Q,R,S,T
# carry the new record's specific data.
$ws.Cells.Item(218, 1).Value = 10
$ws.Cells.Item(218, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(218, 3).Value = "La Araucanía"
$ws.Cells.Item(218, 4).Value = 44510
$ws.Cells.Item(218, 5).Value = 9
$ws.Cells.Item(218, 6).Value = "Fruta"
$ws.Cells.Item(218, 7).Value = 100108
$ws.Cells.Item(218, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(218, 9).Value = 100108005
$ws.Cells.Item(218, 10).Value = "Piña"
$ws.Cells.Item(218, 11).Value = "Caramelo"
$ws.Cells.Item(218, 12).Value = "Segunda"
$ws.Cells.Item(218, 13).Value = 45
$ws.Cells.Item(218, 14).Value = 20000
$ws.Cells.Item(218, 15).Value = 20000
$ws.Cells.Item(218, 16).Value = 20000
$ws.Cells.Item(218, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(218, 18).Value = "Bolivia"
$ws.Cells.Item(218, 19).Value = 1429
$ws.Cells.Item(218, 20).Value = 14
